$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the workbook window down a bit (yWindow 450 -> 900) ---
$win = $wb.Windows.Item(1)
$win.Top = 900

# --- Input value changes; dependent formulas (C/A/E/F columns) recalculate automatically ---
$ws.Range("D9").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("D22").Value = 0

# --- Move the active selection from D19 to D17 ---
$ws.Range("D17").Select()

$wb.Save()
